# B6-PowerPoint.pptx edit
#
# 1) Three tables (on slides 14, 15, 16) get their table style switched
#    from {7C2AF123-264F-48BB-8B8D-ED916158A336} to
#    {FCC819BD-61EC-427D-9CCF-E0BB0BBB9792}.
# 2) The deck's theme colour scheme is swapped from the "Integral /
#    Red Violet" palette to the stock "Office" palette (the
#    before/after pair swap the two theme parts' colour schemes; the
#    font scheme and format scheme are already identical between the
#    two themes, so only the 12 theme colours actually need updating).

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables -----------------------------------
$newStyleId = "{FCC819BD-61EC-427D-9CCF-E0BB0BBB9792}"
$tableSlides = 14, 15, 16

foreach ($idx in $tableSlides) {
    $slide = $p.Slides.Item($idx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2. Swap the theme colours ---------------------------------------
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (the stock "Office"
# colour scheme values), in ThemeColorScheme order.
$officeColors = @(
    0,        # dk1      000000
    16777215, # lt1      FFFFFF
    6968388,  # dk2      44546A
    15132391, # lt2      E7E6E6
    13998939, # accent1  5B9BD5
    3243501,  # accent2  ED7D31
    10855845, # accent3  A5A5A5
    49407,    # accent4  FFC000
    12874308, # accent5  4472C4
    4697456,  # accent6  70AD47
    12673797, # hlink    0563C1
    7491477   # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
